# Remove the trailing "Ver no Jupiter ..." / copyright block from the
# end of the document, along with the blank paragraphs that surround it.
#
# Before:
#   ... LOQ4073: Química Geral II (Requisito fraco)
#   <empty paragraph>
#   Ver no Jupiter Salvar em pdf Salvar em docx
#   © 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution
#   <empty paragraph>
#   <page-break paragraph>
#
# After:
#   ... LOQ4073: Química Geral II (Requisito fraco)
#   <page-break paragraph>

$d = $word.ActiveDocument

$anchorText = "LOQ4073: Química Geral II (Requisito fraco)"
$count = $d.Paragraphs.Count

# Locate the anchor paragraph ("LOQ4073: ...") by scanning with
# indexed access (Paragraph.Next/.Previous are not reliable here).
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
    }
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Find the next paragraph after the anchor that starts a new page
# (the page-break paragraph that should remain immediately after the
# anchor once the block in between is removed).
$stopIndex = -1
for ($i = $anchorIndex + 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.PageBreakBefore) {
        $stopIndex = $i
        break
    }
}

$stopPara = $d.Paragraphs.Item($stopIndex)

$deleteRange = $d.Range($anchorPara.Range.End, $stopPara.Range.Start)
$deleteRange.Delete()
